# Insert a new weekly price record as row 52 ("Ají", Región del Maule,
# $/caja 15 kilos) on sheet "Sheet1". This pushes the former rows 52-61
# down to 53-62, preserving their data untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 52, shifting existing rows 52:61 down to 53:62.
$ws.Rows("52:52").Insert()

# Populate the newly inserted row 52 with the new record.
$ws.Range("A52").Value = 7
$ws.Range("B52").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C52").Value = "Ñuble"
$ws.Range("D52").Value = 44588
$ws.Range("E52").Value = 16
$ws.Range("F52").Value = 100112021
$ws.Range("G52").Value = "Ají"
$ws.Range("H52").Value = "Americana (o)"
$ws.Range("I52").Value = "Primera"
$ws.Range("J52").Value = 60
$ws.Range("K52").Value = 14000
$ws.Range("L52").Value = 14500
$ws.Range("M52").Value = 14250
$ws.Range("N52").Value = "$/caja 15 kilos"
$ws.Range("O52").Value = "Región del Maule"
$ws.Range("P52").Value = 950
$ws.Range("Q52").Value = 15
$ws.Range("R52").Value = "Hortaliza"

# Make sure the date cell keeps the workbook's datetime display format (s=2),
# matching the style used by every other date cell in column D.
$ws.Range("D52").NumberFormat = $ws.Range("D53").NumberFormat
